$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the password value in cell B1
$ws.Range("B1").Value = "password=1"

# Update the selected/active cell in the sheet view
$ws.Range("C6").Select()
